# Shift all timestamps in column A by +3 days (data re-fetched for a later
# date window) and refresh the corresponding "Actual Production (MW)"
# values in column B for the quarter-hourly rows that now have new
# forecasted figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Actual Production (MW) values for rows 2..56 (data rows 1..55 of the
# series). Rows 57..97 stay at 0 (unchanged by the edit).
$newB = @(0,0,0,0,0,0,0,0,0,0,0,0,0,3,11,13,16,15,21,26,27,36,41,47,54,56,55,51,61,70,74,74,58,50,41,34,37,36,38,47,37,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($r = 2; $r -le 97; $r++) {
    # Shift the date/time serial in column A forward by 3 days.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value2 = $cellA.Value2 + 3

    if ($r -le 56) {
        $ws.Cells.Item($r, 2).Value2 = $newB[$r - 2]
    }
}
